$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl = $m.CustomLayouts.Item(2)
$cl.Shapes.Item("Comment").Delete()
$cl.Shapes.Item("Plots").Delete()
$cl.Shapes.Item("MolStructure").Delete()
